$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new work-log entry was recorded for 2014-03-19. It belongs right after the
# existing last data row (row 77), before the blank spacer / summary rows.
# Insert a fresh row at 78, which shifts the spacer + the three summary rows
# (formerly 78-81) down to 79-82, and auto-adjusts the SUM(...) / division
# formulas that referenced them.
$ws.Rows.Item(78).Insert()

# Populate the newly inserted row 78 with the new entry's data.
$ws.Range("A78").Value = 2014
$ws.Range("B78").Value = 3
$ws.Range("C78").Value = 19
$ws.Range("D78").Value = 0.61458333333333337
$ws.Range("E78").Value = 0.75
$ws.Range("F78").Formula = "=(E78-D78)*24*60"
$ws.Range("G78").Formula = "=F78/60"

# Match the author's recorded selection after making the edit.
[void]$ws.Range("F78").Select()
